$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.593.84"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.028.43"
$ws.Range("E3").Value = "  +1.51%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.07"
$ws.Range("E5").Value = "  -8.84%  "

# Row 6
$ws.Range("E6").Value = "  -2.20%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.38"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  -0.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.35"
$ws.Range("E10").Value = "  +3.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  -1.42%  "

# Row 12
$ws.Range("E12").Value = "  -1.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.331.73"
$ws.Range("E13").Value = "  +1.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.32"
$ws.Range("E14").Value = "  +1.60%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.31"
$ws.Range("E15").Value = "  -3.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.763"
$ws.Range("E16").Value = "  -3.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.12"
$ws.Range("E17").Value = "  -0.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.025.34"
$ws.Range("E18").Value = "  +1.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.771.11"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.81"
$ws.Range("E20").Value = "  -3.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.55"
$ws.Range("E21").Value = "  +10.24%  "

# Row 22
$ws.Range("E22").Value = "  -2.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.51"
$ws.Range("E23").Value = "  -5.25%  "

# Row 24
$ws.Range("E24").Value = "  +0.21%  "

# Row 25
$ws.Range("E25").Value = "  +1.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -5.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.64"
$ws.Range("E27").Value = "  -1.15%  "

# Row 28
$ws.Range("E28").Value = "  +8.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.67"
$ws.Range("E29").Value = "  -1.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.99"
$ws.Range("E30").Value = "  -1.97%  "

# Row 31
$ws.Range("E31").Value = "  +1.75%  "

# Row 32
$ws.Range("E32").Value = "  -1.37%  "

# Row 33
$ws.Range("E33").Value = "  -3.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("E34").Value = "  -4.67%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +5.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.28"
$ws.Range("E36").Value = "  -1.82%  "

# Row 37
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  -4.44%  "

# Row 39
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -2.89%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("E40").Value = "  +6.73%  "

# Row 41
$ws.Range("E41").Value = "  -3.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0942"
$ws.Range("E42").Value = "  +3.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.464.47"
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.10"
$ws.Range("E44").Value = "  +6.43%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0204"
$ws.Range("E45").Value = "  -1.79%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.11"
$ws.Range("E46").Value = "  +39.82%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -4.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.72"
$ws.Range("E48").Value = "  +2.15%  "

# Row 49
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("E50").Value = "  -0.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.89"
$ws.Range("E51").Value = "  +0.42%  "
